$d = $word.ActiveDocument

# 1. "Knowledge of many computer languages (" -> "Computer programming expertise ("
$d.Content.Find.Execute("Knowledge of many computer languages (", $true, $false, $false, $false, $false, $true, 1, $false, "Computer programming expertise (", 2)

# 2. Replace the VBA/Python language list (this span crosses the spellStart/gramStart ... spellEnd/gramEnd
#    proofErr markers on both ends, so they are fully removed rather than left dangling)
$d.Content.Find.Execute(", VBA,Python, etc.", $true, $false, $false, $false, $false, $true, 1, $false, ", C, VBA, Python, Matlab, Assembly, …", 2)
